$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Khalili")

# New rows describing search-post / get-post-comments work
$ws.Range("A84").Value = "FragmentSearchResult"
$ws.Range("B84").Value = "lazy load implementation"

$ws.Range("A85").Value = "FragmentSearchResult"
$ws.Range("B85").Value = "list scroll: FORCE CLOSE"

# Scroll / selection state matches the author's saved view
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 79
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("A86").Select()
